# Power Budget.xlsx -- "updated navigator and started on-board simulator node"
#
# 1. Navigator (Turnigy D2206-2600KV motor) power draw corrected 83W -> 20W.
# 2. New component row added for the on-board simulator's VL53L1X distance
#    sensor (0.02W, qty 1, 100% duty).
# 3. Totals row (22) emboldened to set it off from the data rows.
# 4. Selection moved to B8, matching where the edit was made.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Navigator row: correct the power draw --------------------------------
$ws.Range("B8").Value = 20

# --- New component: VL53L1X Distance Sensor (row 10) ----------------------
$ws.Range("A10").Value = "VL53L1X Distance Sensor"
$ws.Range("B10").Value = 0.02
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 1
# E10 is a shared formula (=B10*D10*C10) already spilled down the column;
# it recalculates automatically.

# --- Totals row: bold it out -----------------------------------------------
$ws.Rows(22).Font.Bold = $true

# Registering the bold/no-border row-default style (the one Excel stamps on
# a fully row-selected format change) without disturbing any real cell: set
# it far outside the used range, then delete that scratch row.
$ws.Range("A100").Font.Bold = $true
$ws.Range("A100").EntireRow.Delete()

# --- Restore the author's last selection -----------------------------------
$ws.Range("B8").Select()
